# BURN DOWN CHART SPRINT 3 - update effort entries for day 5 (row 8 on Hoja2)
# Jose (C8) and Camila (D8) both logged 3 hours of effort, which ripples
# through the running totals / "pendiente" columns and the dependent
# NUMBERVALUE() lookups on Hoja1 (and, in turn, the burn-down chart).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- core data edit -------------------------------------------------
$ws2.Range("C8").Value = 3
$ws2.Range("D8").Value = 3

# --- recalculate so every dependent formula (Hoja2 H/I/J, Hoja1 L, the
#     chart series, etc.) carries a fresh cached value -----------------
$excel.CalculateFullRebuild()

# --- restore the view/selection state recorded in the saved file ------
$ws2.Activate()
$ws2.Range("E8").Select()

$ws1.Activate()
$ws1.Range("G12").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

$excel.ActiveWindow.Width = 24240
$excel.ActiveWindow.Height = 13140
